$d = $word.ActiveDocument

# The document currently has a single paragraph (holding the _GoBack
# bookmark) followed by the sectPr. We replace that paragraph's content
# with the full set of new paragraphs (an empty lead-in paragraph, the
# title, three body paragraphs, the bookmark paragraph now wrapped with
# text before/after the bookmark, and a trailing empty paragraph) using
# InsertXML so that the bookmark is preserved exactly where it belongs
# and run/paragraph formatting comes out byte-for-byte as wanted.

$target = $d.Paragraphs(1)

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = @"
<w:p $wns/>
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>A Cat a Parrot and a Bag of Seed</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t xml:space="preserve">Defining the problem. </w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>A man needs to transport 3 items across the river to the other side but the boat only has room for him and 1 item.</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>My insight is he will have to make 3 trips</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t xml:space="preserve"> or find a bigger boat.</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t xml:space="preserve">The overall goal is to get </w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>him and the three items to the other side of the river.</w:t>
  </w:r>
</w:p>
<w:p $wns/>
"@

[void]$target.Range.InsertXML($xml)
